$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Iconos en el menu (no se alinean con el texto)"
#    list item paragraph entirely (it sits right after "General").
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Iconos en el men*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: it currently sits at the end of the
#    "Decidir como mostrar..." paragraph (right before its paragraph
#    mark). It needs to move two paragraphs later, to the start of the
#    very last (empty) paragraph of the document, so that the two
#    trailing empty paragraphs now precede it instead of following it.
#
#    Directly calling Bookmarks.Add on a range inside an empty
#    paragraph can misplace the bookmark, so we first stash a marker
#    run of text at the target location, anchor the bookmark to that
#    non-empty range, and then remove the marker text again - leaving
#    the (now correctly positioned) collapsed bookmark behind.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertBefore("ZZMARKERZZ")

$lastPara = $d.Paragraphs.Last
$markerRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + 10)
$d.Bookmarks.Add("_GoBack", $markerRange)

$lastPara = $d.Paragraphs.Last
$markerRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + 10)
$markerRange.Text = ""
